# MPMC_queue.xlsx - populate Sheet1 with the algorithm description table.
# Cell values are written in "first use" order so that the generated
# shared-strings table lines up with how a human would naturally fill
# the sheet in (content rows first, header row last).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Title
$ws.Range("A1").Value = "Algo:"
$ws.Range("B1").Value = "MultiProducersMultiConsumersUnlimitedLockFreeQueue_v4"

# size = 0 scenario
$ws.Range("A5").Value = "size = 0"
$ws.Range("B5").Value = "first_a = nullptr, last_a = nullptr"

# size = 1 scenario
$ws.Range("A15").Value = "size = 1"
$ws.Range("B15").Value = "first_a = last_a"

# Case labels under each scenario
$ws.Range("A7").Value = "case 1"
$ws.Range("A12").Value = "case 2"

# Header row (row 3) - bold, blue font
$ws.Range("C3").Value = "pop()"
$ws.Range("C3").Font.Bold = $true
$ws.Range("C3").Font.Color = 12611584

$ws.Range("B3").Value = "push()"
$ws.Range("B3").Font.Bold = $true
$ws.Range("B3").Font.Color = 12611584

$ws.Range("D3").Value = "Comments"
$ws.Range("D3").Font.Bold = $true
$ws.Range("D3").Font.Color = 12611584

# Case labels repeated for the size = 1 scenario
$ws.Range("A17").Value = "case 1"
$ws.Range("A19").Value = "case 2"

# Column widths for the three content columns
$ws.Columns("B").ColumnWidth = 73.5
$ws.Columns("C").ColumnWidth = 72.66666666666667
$ws.Columns("D").ColumnWidth = 53.833333333333336

# Selection + print orientation
$ws.Range("B7").Select() | Out-Null
$ws.PageSetup.Orientation = 1
